$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# (matching the original inlineStr/text cell type) by temporarily applying a
# text number format, then reverting the style afterwards so no visible
# formatting change is introduced.
$numericLookingCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D11",
    "D12",
    "D15",
    "D17",
    "D18",
    "D19",
    "D20",
    "D21",
    "D22",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D45",
    "D47",
    "D48",
    "D49",
    "D50"
)
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value2 = "29.738.03"
$ws.Range("E2").Value2 = "  -0.56%  "
$ws.Range("D3").Value2 = "1.885.86"
$ws.Range("E3").Value2 = "  -0.90%  "
$ws.Range("D4").Value2 = "1.000"
$ws.Range("E4").Value2 = "  -0.02%  "
$ws.Range("D5").Value2 = "0.7950"
$ws.Range("E5").Value2 = "  -0.85%  "
$ws.Range("D6").Value2 = "241.54"
$ws.Range("E6").Value2 = "  +0.47%  "
$ws.Range("D7").Value2 = "1.000"
$ws.Range("E7").Value2 = "  -0.01%  "
$ws.Range("E8").Value2 = "  +1.88%  "
$ws.Range("D9").Value2 = "25.52"
$ws.Range("E9").Value2 = "  -2.83%  "
$ws.Range("D10").Value2 = "0.07019"
$ws.Range("E10").Value2 = "  +0.11%  "
$ws.Range("D11").Value2 = "0.08041"
$ws.Range("E11").Value2 = "  +0.62%  "
$ws.Range("D12").Value2 = "0.7684"
$ws.Range("E12").Value2 = "  +4.12%  "
$ws.Range("D13").Value2 = "1.887.44"
$ws.Range("E13").Value2 = "  -0.87%  "
$ws.Range("E14").Value2 = "  +2.70%  "
$ws.Range("D15").Value2 = "91.81"
$ws.Range("E15").Value2 = "  -0.34%  "
$ws.Range("D16").Value2 = "29.748.08"
$ws.Range("E16").Value2 = "  -0.53%  "
$ws.Range("B17").Value2 = "Uniswap"
$ws.Range("C17").Value2 = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value2 = "5.933"
$ws.Range("E17").Value2 = "  +1.67%  "
$ws.Range("B18").Value2 = "Avalanche"
$ws.Range("C18").Value2 = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value2 = "13.77"
$ws.Range("E18").Value2 = "  -1.10%  "
$ws.Range("D19").Value2 = "242.52"
$ws.Range("E19").Value2 = "  -0.74%  "
$ws.Range("D20").Value2 = "0.000007682"
$ws.Range("E20").Value2 = "  -1.17%  "
$ws.Range("B21").Value2 = "Dai"
$ws.Range("C21").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value2 = "1.001"
$ws.Range("E21").Value2 = "  -0.05%  "
$ws.Range("B22").Value2 = "Chainlink"
$ws.Range("C22").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value2 = "8.141"
$ws.Range("E22").Value2 = "  +18.19%  "
$ws.Range("B23").Value2 = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value2 = "2.144.38"
$ws.Range("E23").Value2 = "  -0.81%  "
$ws.Range("D24").Value2 = "1.001"
$ws.Range("E24").Value2 = "  +0.06%  "
$ws.Range("D25").Value2 = "0.1636"
$ws.Range("E25").Value2 = "  +11.89%  "
$ws.Range("D26").Value2 = "9.279"
$ws.Range("E26").Value2 = "  +1.19%  "
$ws.Range("D27").Value2 = "163.76"
$ws.Range("E27").Value2 = "  -2.71%  "
$ws.Range("D28").Value2 = "18.61"
$ws.Range("E28").Value2 = "  -1.16%  "
$ws.Range("D29").Value2 = "2.052"
$ws.Range("E29").Value2 = "  -0.20%  "
$ws.Range("D30").Value2 = "1.368"
$ws.Range("E30").Value2 = "  +0.95%  "
$ws.Range("D31").Value2 = "1.532"
$ws.Range("E31").Value2 = "  +1.40%  "
$ws.Range("D32").Value2 = "4.421"
$ws.Range("E32").Value2 = "  +3.53%  "
$ws.Range("D33").Value2 = "0.05655"
$ws.Range("E33").Value2 = "  +2.63%  "
$ws.Range("D34").Value2 = "4.075"
$ws.Range("E34").Value2 = "  +0.73%  "
$ws.Range("D35").Value2 = "1.259"
$ws.Range("E35").Value2 = "  +0.33%  "
$ws.Range("D36").Value2 = "0.7344"
$ws.Range("E36").Value2 = "  +1.04%  "
$ws.Range("D37").Value2 = "0.9998"
$ws.Range("E37").Value2 = "  +0.04%  "
$ws.Range("D38").Value2 = "2.711"
$ws.Range("E38").Value2 = "  -0.42%  "
$ws.Range("D39").Value2 = "0.01917"
$ws.Range("E39").Value2 = "  +0.15%  "
$ws.Range("D40").Value2 = "2.768"
$ws.Range("E40").Value2 = "  -0.51%  "
$ws.Range("D41").Value2 = "0.4399"
$ws.Range("E41").Value2 = "  +0.17%  "
$ws.Range("D42").Value2 = "71.88"
$ws.Range("E42").Value2 = "  -0.06%  "
$ws.Range("D43").Value2 = "5.820"
$ws.Range("E43").Value2 = "  -2.14%  "
$ws.Range("E44").Value2 = "  +0.03%  "
$ws.Range("D45").Value2 = "0.8383"
$ws.Range("E45").Value2 = "  +0.48%  "
$ws.Range("D46").Value2 = "1.029.55"
$ws.Range("E46").Value2 = "  +5.48%  "
$ws.Range("D47").Value2 = "101.81"
$ws.Range("E47").Value2 = "  +1.02%  "
$ws.Range("D48").Value2 = "1.852"
$ws.Range("E48").Value2 = "  -1.48%  "
$ws.Range("D49").Value2 = "9.857"
$ws.Range("E49").Value2 = "  +1.79%  "
$ws.Range("D50").Value2 = "7.431"
$ws.Range("E50").Value2 = "  -1.31%  "
$ws.Range("D51").Value2 = "2.041.26"
$ws.Range("E51").Value2 = "  -0.80%  "

foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
